$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old third row / second column formatting remnants
$ws.Range("A1:C3").ClearContents()

# Populate in the order that matches the shared-string table build order:
# A1, A2, B1, C1, B2, C2
$ws.Range("A1").Value = "ScenarioName"
$ws.Range("A2").Value = "KindleEBook"
$ws.Range("B1").Value = "data1"
$ws.Range("C1").Value = "data2"
$ws.Range("B2").Value = "value1"
$ws.Range("C2").Value = "value2"

$ws.Range("C2").Select() | Out-Null
